$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164 (1-based), shifting existing rows 164-220 down to 165-221.
$ws.Rows.Item(164).Insert()

# Copy the data that was in the old row 164 (now shifted to row 165) into the new row 164,
# except columns D (Fecha) and J (Volumen) which get new values.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(164, $col).Value2 = $ws.Cells.Item(165, $col).Value2
}

$ws.Cells.Item(164, 4).Value2 = 44900
$ws.Cells.Item(164, 10).Value2 = 15
